# Edit: split the single "Annotation" column into a 3-level annotation
# hierarchy (ann_level_3, ann_level_2, ann_level_1) by inserting a new
# column C (a duplicate of the existing Annotation column B), then
# renaming the header cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C. This shifts the old "Broad"
# column (C) and everything to its right one column to the right, and
# copies the formatting of column B into the newly created column C.
$ws.Columns.Item(3).Insert()

# The insert operation duplicated column B's formatting into the new
# column C; now copy the actual values of column B into column C,
# row by row, so both hold the same annotation text (rows 2-22).
# (Value2 is used because bulk Range-to-Range .Value array assignment
# is not reliable in this runtime.)
for ($r = 2; $r -le 22; $r++) {
    $ws.Cells.Item($r, 3).Value2 = $ws.Cells.Item($r, 2).Value2
}

# Rename the header row to reflect the new multi-level annotation
# scheme.
$ws.Range("B2").Value2 = "ann_level_3"
$ws.Range("C2").Value2 = "ann_level_2"
$ws.Range("D2").Value2 = "ann_level_1"

# Update selection / active cell to match the saved view state.
$ws.Range("C4").Select()
